$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column B for the new rows as Text so that numeric-looking
# account numbers are preserved as strings (matches existing column B data).
$ws.Range("B245:B264").NumberFormat = "@"

$ws.Cells.Item(245, 1).Value = '2026-02-16 14:26:40'
$ws.Cells.Item(245, 2).Value = '237675239360'
$ws.Cells.Item(245, 3).Value = 'EDITH LAURE MATCHINDA NGUEMETA'
$ws.Cells.Item(245, 4).Value = 201577

$ws.Cells.Item(246, 1).Value = '2026-02-16 15:55:11'
$ws.Cells.Item(246, 2).Value = '237675396752'
$ws.Cells.Item(246, 3).Value = 'BENEDICTE CHANTAL MANTSANG'
$ws.Cells.Item(246, 4).Value = 69826

$ws.Cells.Item(247, 1).Value = '2026-02-16 15:18:36'
$ws.Cells.Item(247, 2).Value = '237675626141'
$ws.Cells.Item(247, 3).Value = 'FLORENCE NGUEFACK'
$ws.Cells.Item(247, 4).Value = 15874

$ws.Cells.Item(248, 1).Value = '2026-02-17 00:14:35'
$ws.Cells.Item(248, 2).Value = '237676840777'
$ws.Cells.Item(248, 3).Value = 'ETP109 ETP'
$ws.Cells.Item(248, 4).Value = 1216408

$ws.Cells.Item(249, 1).Value = '2026-02-16 15:01:29'
$ws.Cells.Item(249, 2).Value = '237677833877'
$ws.Cells.Item(249, 3).Value = 'MEDJOM TAGNE MICHELLE GUILENE LA NEGRESSE SARL'
$ws.Cells.Item(249, 4).Value = 221

$ws.Cells.Item(250, 1).Value = '2026-02-16 12:48:07'
$ws.Cells.Item(250, 2).Value = '237678854978'
$ws.Cells.Item(250, 3).Value = 'ELVIS FEUDJIO'
$ws.Cells.Item(250, 4).Value = 399496

$ws.Cells.Item(251, 1).Value = '2026-02-16 14:57:14'
$ws.Cells.Item(251, 2).Value = '237679422591'
$ws.Cells.Item(251, 3).Value = 'ETS LE CONTENT 42'
$ws.Cells.Item(251, 4).Value = 403900

$ws.Cells.Item(252, 1).Value = '2026-02-16 14:11:15'
$ws.Cells.Item(252, 2).Value = '237650353920'
$ws.Cells.Item(252, 3).Value = 'MENIAPI HELENE EDOSSINE TOP MOBIL TELECOM'
$ws.Cells.Item(252, 4).Value = 1201274

$ws.Cells.Item(253, 1).Value = '2026-02-16 22:48:51'
$ws.Cells.Item(253, 2).Value = '237651927448'
$ws.Cells.Item(253, 3).Value = 'charity aben awalah'
$ws.Cells.Item(253, 4).Value = 306869

$ws.Cells.Item(254, 1).Value = '2026-02-16 12:11:27'
$ws.Cells.Item(254, 2).Value = '237653294562'
$ws.Cells.Item(254, 3).Value = 'NANHOU KEMAYOU AVIGAEL ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Cells.Item(254, 4).Value = 1000736

$ws.Cells.Item(255, 1).Value = '2026-02-16 13:46:57'
$ws.Cells.Item(255, 2).Value = '237678046498'
$ws.Cells.Item(255, 3).Value = 'MFS SOCAVER'
$ws.Cells.Item(255, 4).Value = 50134

$ws.Cells.Item(256, 1).Value = '2026-02-16 14:57:54'
$ws.Cells.Item(256, 2).Value = '237679428698'
$ws.Cells.Item(256, 3).Value = 'ETS LE CONTENT 29'
$ws.Cells.Item(256, 4).Value = 32

$ws.Cells.Item(257, 1).Value = '2026-02-16 15:10:13'
$ws.Cells.Item(257, 2).Value = '237679551262'
$ws.Cells.Item(257, 3).Value = 'LA NEGRESSE LTDLA CBOX R1 MEGNE JUDITH'
$ws.Cells.Item(257, 4).Value = 67562

$ws.Cells.Item(258, 1).Value = '2026-02-16 11:15:26'
$ws.Cells.Item(258, 2).Value = '237680574202'
$ws.Cells.Item(258, 3).Value = 'TOUMEWO SAMUEL'
$ws.Cells.Item(258, 4).Value = 308260

$ws.Cells.Item(259, 1).Value = '2026-02-16 13:33:54'
$ws.Cells.Item(259, 2).Value = '237681118330'
$ws.Cells.Item(259, 3).Value = 'SAHA NDESA JONAS LTDLA_POLAS_OTH_NDOGBONG SERIE'
$ws.Cells.Item(259, 4).Value = 157629

$ws.Cells.Item(260, 1).Value = '2026-02-16 16:05:13'
$ws.Cells.Item(260, 2).Value = '237674446293'
$ws.Cells.Item(260, 3).Value = 'jean michel mba'
$ws.Cells.Item(260, 4).Value = 31986

$ws.Cells.Item(261, 1).Value = '2026-02-16 14:56:55'
$ws.Cells.Item(261, 2).Value = '237679085953'
$ws.Cells.Item(261, 3).Value = 'MADELEINE NKOUADJIO'
$ws.Cells.Item(261, 4).Value = 24046

$ws.Cells.Item(262, 1).Value = '2026-02-16 16:39:04'
$ws.Cells.Item(262, 2).Value = '237681662761'
$ws.Cells.Item(262, 3).Value = 'LUC BAYOMOCK'
$ws.Cells.Item(262, 4).Value = 67626

$ws.Cells.Item(263, 1).Value = '2026-02-16 17:06:47'
$ws.Cells.Item(263, 2).Value = '237682975726'
$ws.Cells.Item(263, 3).Value = 'LA NEGRESSE SARL NYOUNG JOSEPH CLOTAIRE'
$ws.Cells.Item(263, 4).Value = 17785

$ws.Cells.Item(264, 1).Value = '2026-02-16 10:41:59'
$ws.Cells.Item(264, 2).Value = '237683075075'
$ws.Cells.Item(264, 3).Value = 'GORBATCHEV NGUETSA KOUAKAM'
$ws.Cells.Item(264, 4).Value = 13193

